$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns R1 (Coupons) and S1 (Discounts)
$ws.Range("R1").Value = "Coupons"
$ws.Range("S1").Value = "Discounts"

# New row 14: "PLU Lookup" test case
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "PLU Lookup"
$ws.Range("C14").Value = "Dine In"
$ws.Range("D14").Value = "Coupons Discounts, Promo Lookup"
$ws.Range("E14").Value = "Classic Sandwch, Spicy Sandwch, Classic Sandwch, Spicy Sandwich, 4 Biscuit, Cajun Rice Lg"
$ws.Range("G14").Value = "D-Sandwich Family Meal"
$ws.Range("K14").Value = "D-Sandwich Family Meal"
$ws.Range("L14").Value = "Classic Sandwch, Spicy Sandwch"
$ws.Range("N14").Value = "Classic Sandwch"
$ws.Range("P14").Value = "YES"
$ws.Range("Q14").Value = "NO"
$ws.Range("R14").Value = 1471

# New row 15: "X% Off" test case
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "X% Off"
$ws.Range("C15").Value = "Dine In"
$ws.Range("D15").Value = "Sandwiches Seafood"
$ws.Range("E15").Value = "CLASSIC CHICKEN SANDWICH, Crispy, Classic Chicken Sandwich"
$ws.Range("F15").Value = "Sandwich Only"
$ws.Range("G15").Value = "Classic Sandwch"
$ws.Range("K15").Value = "Classic Sandwch"
$ws.Range("L15").Value = "Classic Sandwch"
$ws.Range("M15").Value = "Classic Sandwich"
$ws.Range("N15").Value = "Classic Sandwch"
$ws.Range("P15").Value = "YES"
$ws.Range("Q15").Value = "NO"
$ws.Range("S15").Value = "25% Discount"

# Match the final selection state recorded in the saved workbook
$ws.Range("Q15").Select()
